# Merge the "formativas." runs back into the main paragraph run on the
# "Objetivo General" slide (slide 5), removing the underline formatting
# that set "formativas" apart and the trailing "." run, so the whole
# sentence becomes a single, uniformly formatted run.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$target = "Diseñar e implementar un sistema digital para la gestión de minutas, el control de asistencia y el registro de incidentes en los ambientes del SENA, con el fin de optimizar la organización institucional, mejorar la comunicación entre actores y facilitar el seguimiento eficiente de las actividades formativas."

# First collapse the whole paragraph down to a single run (this adopts the
# formatting of the first existing run), then restore the full text. Since
# there is now only one run to diff against, the final assignment keeps
# everything as one run instead of re-splitting it around "formativas".
$tr.Text = "X"
$tr.Text = $target
